$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ C = 117868605.5867971; D = 1.832837670377314;  E = 4298.174879869647; F = 4298.174879869647 }
    3  = @{ C = 117853749.5525672; D = 1.672161674023858;  E = 4711.77518178011;  F = 9009.950061649757 }
    4  = @{ C = 117601884.879902;  D = 1.521127160740514;  E = 4985.822892365274; F = 13995.77295401503 }
    5  = @{ C = 117822492.0244499; D = 1.634850098626581;  E = 4630.314699945834; F = 18626.08765396087 }
    6  = @{ C = 117790726.6878049; D = 1.633760465320453;  E = 5365.621325087924; F = 23991.70897904879 }
    7  = @{ C = 117776698.0440098; D = 1.470056583249274;  E = 5963.841288213086; F = 29955.55026726187 }
    8  = @{ C = 117826938.0808824; D = 1.676471650687426;  E = 4160.954526898982; F = 34116.50479416086 }
    9  = @{ C = 117601998.8365854; D = 1.669251038708748;  E = 3983.072277252608; F = 38099.57707141346 }
    10 = @{ C = 117994946.7188264; D = 1.432430522853284;  E = 4705.922704065689; F = 42805.49977547915 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
}
